$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking correct count (row 11, column B)
$ws.Range("B11").Value = 5

# Update total marks (row 12, column B) and the corr/total fraction text (E12)
$ws.Range("B12").Value = 120
$ws.Range("E12").Value = "120/140"
